$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (D value, E value) mapping derived from the target diff (rows 33-94)
$values = @{
    33 = @(1, 1)
    34 = @(0, 1)
    35 = @(0, 1)
    36 = @(0, 1)
    37 = @(0, 1)
    38 = @(1, 1)
    39 = @(0, 0)
    40 = @(1, 0)
    41 = @(1, 1)
    42 = @(0, 1)
    43 = @(0, 0)
    44 = @(1, 1)
    45 = @(0, 0)
    46 = @(0, 1)
    47 = @(1, 0)
    48 = @(1, 1)
    49 = @(1, 1)
    50 = @(1, 0)
    51 = @(1, 1)
    52 = @(0, 0)
    53 = @(0, 0)
    54 = @(1, 1)
    55 = @(0, 0)
    56 = @(0, 0)
    57 = @(0, 1)
    58 = @(0, 0)
    59 = @(1, 1)
    60 = @(1, 0)
    61 = @(0, 1)
    62 = @(1, 0)
    63 = @(1, 1)
    64 = @(1, 0)
    65 = @(0, 1)
    66 = @(0, 0)
    67 = @(2, 0)
    68 = @(1, 0)
    69 = @(1, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(0, 1)
    73 = @(1, 0)
    74 = @(1, 1)
    75 = @(0, 0)
    76 = @(0, 0)
    77 = @(1, 0)
    78 = @(0, 1)
    79 = @(1, 0)
    80 = @(1, 0)
    81 = @(0, 0)
    82 = @(0, 1)
    83 = @(0, 1)
    84 = @(0, 1)
    85 = @(1, 1)
    86 = @(0, 0)
    87 = @(1, 0)
    88 = @(0, 0)
    89 = @(1, 1)
    90 = @(1, 1)
    91 = @(1, 1)
    92 = @(1, 0)
    93 = @(1, 1)
    94 = @(0, 1)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item([int]$row, 4).Value = $pair[0]  # Column D
    $ws.Cells.Item([int]$row, 5).Value = $pair[1]  # Column E
}
